$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 219, shifting the existing rows 219:284 down to 220:285.
$ws.Rows("219:219").Insert()

# Populate the newly inserted row 219 with the new weekly price record.
$ws.Range("A219").Value = 8
$ws.Range("B219").Value = "Terminal La Palmera de La Serena"
$ws.Range("C219").Value = "Coquimbo"
$ws.Range("D219").Value = 44841
$ws.Range("E219").Value = 4
$ws.Range("F219").Value = 100112031
$ws.Range("G219").Value = "Poroto verde"
$ws.Range("H219").Value = "Magnum"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 600
$ws.Range("K219").Value = 34000
$ws.Range("L219").Value = 35000
$ws.Range("M219").Value = 34500
$ws.Range("N219").Value = "$/malla 25 kilos"
$ws.Range("O219").Value = "Perú"
$ws.Range("P219").Value = 1380
$ws.Range("Q219").Value = 25
$ws.Range("R219").Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Range("D219").NumberFormat = $ws.Range("D220").NumberFormat
